# Auto-generated edit script: apply market-price refresh values to Aegis_Profits sheets
# Each row below corresponds to a single leve/recipe whose cached market data (columns H-N)
# was refreshed by the scheduled runner. Values are written as literal numbers; a few rows
# gain or lose a trailing cell (HQ/NQ profit column) because the source no longer produced it.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 169.09091
$ws.Range("I9").Value = 153.33333
$ws.Range("J9").Value = 240
$ws.Range("K9").Value = 153.33333
$ws.Range("L9").Value = 240
$ws.Range("M9").Value = 15.66667000000001
$ws.Range("N9").Value = -578
# Row 86
$ws.Range("H86").Value = 9093994
$ws.Range("I86").Value = 1434.8572
$ws.Range("J86").Value = 25005972
$ws.Range("K86").Value = 1434.8572
$ws.Range("L86").Value = 25005972
$ws.Range("M86").Value = -311.8571999999999
$ws.Range("N86").Value = -25008218
# Row 89
$ws.Range("H89").Value = 9093994
$ws.Range("I89").Value = 1434.8572
$ws.Range("J89").Value = 25005972
$ws.Range("K89").Value = 7174.286
$ws.Range("L89").Value = 125029860
$ws.Range("M89").Value = -1558.286
$ws.Range("N89").Value = -125041092
# Row 98
$ws.Range("H98").Value = 2258.3845
$ws.Range("I98").Value = 1215.3636
$ws.Range("K98").Value = 1215.3636
$ws.Range("M98").Value = 282.6364000000001
# Row 122
$ws.Range("H122").Value = 2258.3845
$ws.Range("I122").Value = 1215.3636
$ws.Range("K122").Value = 3646.0908
$ws.Range("M122").Value = -1196.0908
# Row 129
$ws.Range("H129").Value = 4619.2705
$ws.Range("J129").Value = 2665.5417
$ws.Range("L129").Value = 7996.625100000001
$ws.Range("N129").Value = -17996.6251
# Row 141
$ws.Range("H141").Value = 2114.375
$ws.Range("I141").Value = 1968.75
$ws.Range("K141").Value = 5906.25
$ws.Range("M141").Value = -726.25

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Range("H23").Value = 48503.125
$ws.Range("J23").Value = 29601.4
$ws.Range("L23").Value = 29601.4
$ws.Range("N23").Value = -30119.4
# Row 32
$ws.Range("H32").Value = 3432.7576
$ws.Range("I32").Value = 2875.3223
$ws.Range("J32").Value = 9007.111000000001
$ws.Range("K32").Value = 2875.3223
$ws.Range("L32").Value = 9007.111000000001
$ws.Range("M32").Value = -2588.3223
$ws.Range("N32").Value = -9581.111000000001
# Row 61
$ws.Range("H61").Value = 2310.9487
$ws.Range("I61").Value = 1376.4375
$ws.Range("J61").Value = 2961.0435
$ws.Range("K61").Value = 1376.4375
$ws.Range("L61").Value = 2961.0435
$ws.Range("M61").Value = -1164.4375
$ws.Range("N61").Value = -3385.0435
# Row 74
$ws.Range("H74").Value = 847.28
$ws.Range("I74").Value = 850.95654
$ws.Range("J74").Value = 805
$ws.Range("K74").Value = 850.95654
$ws.Range("L74").Value = 805
$ws.Range("M74").Value = 23.04345999999998
$ws.Range("N74").Value = -2553
# Row 77
$ws.Range("H77").Value = 847.28
$ws.Range("I77").Value = 850.95654
$ws.Range("J77").Value = 805
$ws.Range("K77").Value = 4254.7827
$ws.Range("L77").Value = 4025
$ws.Range("M77").Value = 113.2173000000003
$ws.Range("N77").Value = -12761
# Row 136
$ws.Range("H136").Value = 2310.9487
$ws.Range("I136").Value = 1376.4375
$ws.Range("J136").Value = 2961.0435
$ws.Range("K136").Value = 4129.3125
$ws.Range("L136").Value = 8883.130500000001
$ws.Range("M136").Value = -1579.3125
$ws.Range("N136").Value = -13983.1305

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 5834.613
$ws.Range("I134").Value = 5766.9546
$ws.Range("J134").Value = 6000
$ws.Range("K134").Value = 17300.8638
$ws.Range("L134").Value = 18000
$ws.Range("M134").Value = -14765.8638
$ws.Range("N134").Value = -23070

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 50004590
$ws.Range("I132").Value = 47623652
$ws.Range("J132").Value = 55560110
$ws.Range("K132").Value = 142870956
$ws.Range("L132").Value = 166680330
$ws.Range("M132").Value = -142868426
$ws.Range("N132").Value = -166685390
# Row 134
$ws.Range("H134").Value = 1169.973
$ws.Range("I134").Value = 1051.4166
$ws.Range("J134").Value = 1388.8462
$ws.Range("K134").Value = 3154.2498
$ws.Range("L134").Value = 4166.5386
$ws.Range("M134").Value = -619.2498000000001
$ws.Range("N134").Value = -9236.5386

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1545.3636
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 6000
$ws.Range("N4").Value = -6224
# Row 12
$ws.Range("H12").Value = 35.454544
$ws.Range("I12").Value = 16.5
$ws.Range("J12").Value = 46.285713
$ws.Range("K12").Value = 49.5
$ws.Range("L12").Value = 138.857139
$ws.Range("M12").Value = 123.5
$ws.Range("N12").Value = -484.857139
# Row 80
$ws.Range("H80").Value = 13619
$ws.Range("J80").Value = 13619
$ws.Range("L80").Value = 40857
$ws.Range("N80").Value = -42729
# Row 83
$ws.Range("H83").Value = 13619
$ws.Range("J83").Value = 13619
$ws.Range("L83").Value = 122571
$ws.Range("N83").Value = -131931
# Row 110
$ws.Range("H110").Value = 862.3333
$ws.Range("I110").Value = 862.3333
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2586.9999
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1503.0001
$ws.Range("N110").ClearContents()
# Row 112
$ws.Range("H112").Value = 102071.7
$ws.Range("J112").Value = 2498.5715
$ws.Range("L112").Value = 7495.7145
$ws.Range("N112").Value = -9711.7145
# Row 124
$ws.Range("H124").Value = 2015.3846
$ws.Range("I124").Value = 916.6667
$ws.Range("J124").Value = 4487.5
$ws.Range("K124").Value = 2750.0001
$ws.Range("L124").Value = 13462.5
$ws.Range("M124").Value = 2159.9999
$ws.Range("N124").Value = -23282.5
# Row 131
$ws.Range("H131").Value = 1540.5698
$ws.Range("I131").Value = 2033.3334
$ws.Range("J131").Value = 1522.759
$ws.Range("K131").Value = 6100.0002
$ws.Range("L131").Value = 4568.277
$ws.Range("M131").Value = -1060.0002
$ws.Range("N131").Value = -14648.277
# Row 134
$ws.Range("H134").Value = 3033.258
$ws.Range("I134").Value = 2082.7334
$ws.Range("J134").Value = 3924.375
$ws.Range("K134").Value = 6248.2002
$ws.Range("L134").Value = 11773.125
$ws.Range("M134").Value = -1178.2002
$ws.Range("N134").Value = -21913.125
# Row 136
$ws.Range("H136").Value = 1564.6
$ws.Range("I136").Value = 1349.5555
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4048.6665
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = 1051.3335
$ws.Range("N136").Value = -20700

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 4598.8184
$ws.Range("I132").Value = 4499.5
$ws.Range("J132").Value = 4772.625
$ws.Range("K132").Value = 13498.5
$ws.Range("L132").Value = 14317.875
$ws.Range("M132").Value = -10968.5
$ws.Range("N132").Value = -19377.875

$ws = $wb.Worksheets.Item("LTW")
# Row 137
$ws.Range("H137").Value = 61857.25
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 368.33334
$ws.Range("I23").Value = 368.33334
$ws.Range("K23").Value = 368.33334
$ws.Range("M23").Value = -139.33334
# Row 124
$ws.Range("H124").Value = 43000
$ws.Range("J124").Value = 43000
$ws.Range("L124").Value = 43000
$ws.Range("N124").Value = -52820
# Row 136
$ws.Range("H136").Value = 712.5625
$ws.Range("I136").Value = 417.45715
$ws.Range("K136").Value = 1252.37145
$ws.Range("M136").Value = 1297.62855
